$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Append a new row of metric data at row 72
$ws.Range("A72").Value = "2025-04-29 10:49:31"
$ws.Range("B72").Value = 205
